$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = 0.953125
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = 0.90625
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = 0.890625
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = 0.84375
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = 0.859375
$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = 0.8125
$ws.Cells.Item(10, 1).Value = 6
$ws.Cells.Item(10, 2).Value = 0.8125
$ws.Cells.Item(11, 1).Value = 7
$ws.Cells.Item(11, 2).Value = 0.78125
$ws.Cells.Item(12, 1).Value = 8
$ws.Cells.Item(12, 2).Value = 0.75
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = 0.75
$ws.Cells.Item(14, 1).Value = 10
$ws.Cells.Item(14, 2).Value = 0.765625
$ws.Cells.Item(15, 1).Value = 11
$ws.Cells.Item(15, 2).Value = 0.71875
$ws.Cells.Item(16, 1).Value = 12
$ws.Cells.Item(16, 2).Value = 0.71875
$ws.Cells.Item(17, 1).Value = 13
$ws.Cells.Item(17, 2).Value = 0.75
$ws.Cells.Item(18, 1).Value = 14
$ws.Cells.Item(18, 2).Value = 0.75
$ws.Cells.Item(19, 1).Value = 15
$ws.Cells.Item(19, 2).Value = 0.640625
$ws.Cells.Item(20, 1).Value = 16
$ws.Cells.Item(20, 2).Value = 0.703125
$ws.Cells.Item(21, 1).Value = 17
$ws.Cells.Item(21, 2).Value = 0.671875
$ws.Cells.Item(22, 1).Value = 18
$ws.Cells.Item(22, 2).Value = 0.640625
$ws.Cells.Item(23, 1).Value = 19
$ws.Cells.Item(23, 2).Value = 0.640625
$ws.Cells.Item(24, 1).Value = 20
$ws.Cells.Item(24, 2).Value = 0.59375
$ws.Cells.Item(25, 1).Value = 21
$ws.Cells.Item(25, 2).Value = 0.59375
$ws.Cells.Item(26, 1).Value = 22
$ws.Cells.Item(26, 2).Value = 0.5625
$ws.Cells.Item(27, 1).Value = 23
$ws.Cells.Item(27, 2).Value = 0.625
$ws.Cells.Item(28, 1).Value = 24
$ws.Cells.Item(28, 2).Value = 0.640625
$ws.Cells.Item(29, 1).Value = 25
$ws.Cells.Item(29, 2).Value = 0.640625
$ws.Cells.Item(30, 1).Value = 26
$ws.Cells.Item(30, 2).Value = 0.625
$ws.Cells.Item(31, 1).Value = 27
$ws.Cells.Item(31, 2).Value = 0.640625
$ws.Cells.Item(32, 1).Value = 28
$ws.Cells.Item(32, 2).Value = 0.640625
$ws.Cells.Item(33, 1).Value = 29
$ws.Cells.Item(33, 2).Value = 0.640625
$ws.Cells.Item(34, 1).Value = 30
$ws.Cells.Item(34, 2).Value = 0.640625
$ws.Cells.Item(35, 1).Value = 31
$ws.Cells.Item(35, 2).Value = 0.625
$ws.Cells.Item(36, 1).Value = 32
$ws.Cells.Item(36, 2).Value = 0.625
$ws.Cells.Item(37, 1).Value = 33
$ws.Cells.Item(37, 2).Value = 0.625
$ws.Cells.Item(38, 1).Value = 34
$ws.Cells.Item(38, 2).Value = 0.609375
$ws.Cells.Item(39, 1).Value = 35
$ws.Cells.Item(39, 2).Value = 0.609375
$ws.Cells.Item(40, 1).Value = 36
$ws.Cells.Item(40, 2).Value = 0.609375
$ws.Cells.Item(41, 1).Value = 37
$ws.Cells.Item(41, 2).Value = 0.609375
$ws.Cells.Item(42, 1).Value = 38
$ws.Cells.Item(42, 2).Value = 0.609375
$ws.Cells.Item(43, 1).Value = 39
$ws.Cells.Item(43, 2).Value = 0.609375
$ws.Cells.Item(44, 1).Value = 40
$ws.Cells.Item(44, 2).Value = 0.609375
$ws.Cells.Item(45, 1).Value = 41
$ws.Cells.Item(45, 2).Value = 0.609375
$ws.Cells.Item(46, 1).Value = 42
$ws.Cells.Item(46, 2).Value = 0.625
$ws.Cells.Item(47, 1).Value = 43
$ws.Cells.Item(47, 2).Value = 0.625
$ws.Cells.Item(48, 1).Value = 44
$ws.Cells.Item(48, 2).Value = 0.625
$ws.Cells.Item(49, 1).Value = 45
$ws.Cells.Item(49, 2).Value = 0.625
$ws.Cells.Item(50, 1).Value = 46
$ws.Cells.Item(50, 2).Value = 0.625
$ws.Cells.Item(51, 1).Value = 47
$ws.Cells.Item(51, 2).Value = 0.625
$ws.Cells.Item(52, 1).Value = 48
$ws.Cells.Item(52, 2).Value = 0.625
$ws.Cells.Item(53, 1).Value = 49
$ws.Cells.Item(53, 2).Value = 0.625
$ws.Cells.Item(54, 1).Value = 50
$ws.Cells.Item(54, 2).Value = 0.625
$ws.Cells.Item(55, 1).Value = 51
$ws.Cells.Item(55, 2).Value = 0.609375
$ws.Cells.Item(56, 1).Value = 52
$ws.Cells.Item(56, 2).Value = 0.609375
$ws.Cells.Item(57, 1).Value = 53
$ws.Cells.Item(57, 2).Value = 0.609375
$ws.Cells.Item(58, 1).Value = 54
$ws.Cells.Item(58, 2).Value = 0.609375
$ws.Cells.Item(59, 1).Value = 55
$ws.Cells.Item(59, 2).Value = 0.625
$ws.Cells.Item(60, 1).Value = 56
$ws.Cells.Item(60, 2).Value = 0.625
$ws.Cells.Item(61, 1).Value = 57
$ws.Cells.Item(61, 2).Value = 0.625
$ws.Cells.Item(62, 1).Value = 58
$ws.Cells.Item(62, 2).Value = 0.625
$ws.Cells.Item(63, 1).Value = 59
$ws.Cells.Item(63, 2).Value = 0.625
$ws.Cells.Item(64, 1).Value = 60
$ws.Cells.Item(64, 2).Value = 0.625
$ws.Cells.Item(65, 1).Value = 61
$ws.Cells.Item(65, 2).Value = 0.625
$ws.Cells.Item(66, 1).Value = 62
$ws.Cells.Item(66, 2).Value = 0.625
$ws.Cells.Item(67, 1).Value = 63
$ws.Cells.Item(67, 2).Value = 0.625
$ws.Cells.Item(68, 1).Value = 64
$ws.Cells.Item(68, 2).Value = 0.625
$ws.Cells.Item(69, 1).Value = 65
$ws.Cells.Item(69, 2).Value = 0.625
$ws.Cells.Item(70, 1).Value = 66
$ws.Cells.Item(70, 2).Value = 0.625
$ws.Cells.Item(71, 1).Value = 67
$ws.Cells.Item(71, 2).Value = 0.625
$ws.Cells.Item(72, 1).Value = 68
$ws.Cells.Item(72, 2).Value = 0.625
$ws.Cells.Item(73, 1).Value = 69
$ws.Cells.Item(73, 2).Value = 0.625
$ws.Cells.Item(74, 1).Value = 70
$ws.Cells.Item(74, 2).Value = 0.625
$ws.Cells.Item(75, 1).Value = 71
$ws.Cells.Item(75, 2).Value = 0.625
$ws.Cells.Item(76, 1).Value = 72
$ws.Cells.Item(76, 2).Value = 0.625
$ws.Cells.Item(77, 1).Value = 73
$ws.Cells.Item(77, 2).Value = 0.625
$ws.Cells.Item(78, 1).Value = 74
$ws.Cells.Item(78, 2).Value = 0.609375
$ws.Cells.Item(79, 1).Value = 75
$ws.Cells.Item(79, 2).Value = 0.609375
$ws.Cells.Item(80, 1).Value = 76
$ws.Cells.Item(80, 2).Value = 0.609375
$ws.Cells.Item(81, 1).Value = 77
$ws.Cells.Item(81, 2).Value = 0.609375
$ws.Cells.Item(82, 1).Value = 78
$ws.Cells.Item(82, 2).Value = 0.609375
$ws.Cells.Item(83, 1).Value = 79
$ws.Cells.Item(83, 2).Value = 0.609375
$ws.Cells.Item(84, 1).Value = 80
$ws.Cells.Item(84, 2).Value = 0.625
$ws.Cells.Item(85, 1).Value = 81
$ws.Cells.Item(85, 2).Value = 0.625
$ws.Cells.Item(86, 1).Value = 82
$ws.Cells.Item(86, 2).Value = 0.625
$ws.Cells.Item(87, 1).Value = 83
$ws.Cells.Item(87, 2).Value = 0.625
$ws.Cells.Item(88, 1).Value = 84
$ws.Cells.Item(88, 2).Value = 0.625
$ws.Cells.Item(89, 1).Value = 85
$ws.Cells.Item(89, 2).Value = 0.625
$ws.Cells.Item(90, 1).Value = 86
$ws.Cells.Item(90, 2).Value = 0.625
$ws.Cells.Item(91, 1).Value = 87
$ws.Cells.Item(91, 2).Value = 0.625
$ws.Cells.Item(92, 1).Value = 88
$ws.Cells.Item(92, 2).Value = 0.625
$ws.Cells.Item(93, 1).Value = 89
$ws.Cells.Item(93, 2).Value = 0.625
$ws.Cells.Item(94, 1).Value = 90
$ws.Cells.Item(94, 2).Value = 0.625
$ws.Cells.Item(95, 1).Value = 91
$ws.Cells.Item(95, 2).Value = 0.625
$ws.Cells.Item(96, 1).Value = 92
$ws.Cells.Item(96, 2).Value = 0.625
$ws.Cells.Item(97, 1).Value = 93
$ws.Cells.Item(97, 2).Value = 0.625
$ws.Cells.Item(98, 1).Value = 94
$ws.Cells.Item(98, 2).Value = 0.625
$ws.Cells.Item(99, 1).Value = 95
$ws.Cells.Item(99, 2).Value = 0.625
$ws.Cells.Item(100, 1).Value = 96
$ws.Cells.Item(100, 2).Value = 0.625
$ws.Cells.Item(101, 1).Value = 97
$ws.Cells.Item(101, 2).Value = 0.625
$ws.Cells.Item(102, 1).Value = 98
$ws.Cells.Item(102, 2).Value = 0.625
$ws.Cells.Item(103, 1).Value = 99
$ws.Cells.Item(103, 2).Value = 0.625
$ws.Cells.Item(104, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(104, 2).Value = 0.625
$ws.Cells.Item(105, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(105, 2).Value = 0.703125
$ws.Cells.Item(106, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(106, 2).Value = 0.65625
$ws.Cells.Item(107, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(107, 2).Value = 0.59375
$ws.Cells.Item(108, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(108, 2).Value = 0.484375
$ws.Cells.Item(109, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(109, 2).Value = 0.65625
$ws.Cells.Item(110, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(110, 2).Value = 0.609375
$ws.Cells.Item(111, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(111, 2).Value = 0.65625
$ws.Cells.Item(112, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(112, 2).Value = 0.609375
$ws.Cells.Item(113, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(113, 2).Value = 0.640625
$ws.Cells.Item(114, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(114, 2).Value = 0.625
$ws.Cells.Item(115, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(115, 2).Value = 0.65625
$ws.Cells.Item(116, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(116, 2).Value = 0.59375
$ws.Cells.Item(117, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(117, 2).Value = 0.625
$ws.Cells.Item(118, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(118, 2).Value = 0.671875
$ws.Cells.Item(119, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(119, 2).Value = 0.640625
$ws.Cells.Item(120, 1).Value = '<__main__.DisplayOutputs object at 0x7fcfe04bb910>'
$ws.Cells.Item(120, 2).Value = 0.5081967213114754

$ws.Cells.Select()
$ws.Cells.Item(16, 9).Activate()
